# CIERRE 10 SEPT 22
# Advance the payroll week header from "SEMANA 35 ... 29 Al 04 DE SEPTIEMBRE 2022"
# to "SEMANA 36 ... 05 Al 11 DE SEPTIEMBRE 2022" and bump the week's extra-hours
# total (K21) from 560 to 980. Cells H9/B27/H27/B43 and the K24 subtotal are
# formulas that recompute automatically from B9/K21.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("recibos")

# Update the week-range label (B9). Downstream formula cells (H9 = B9,
# B27 = B9, H27 = B27, B43 = H27) pick this up automatically on recalc.
$ws.Range("B9").Value = "SEMANA  36  DEL    05      Al   11   DE  SEPTIEMBRE          2022"

# Update the week's extra-hours amount; K24 = SUM(K21:K23) recalculates too.
$ws.Range("K21").Value = 980

# Restore the view position (scrolled down a bit further, new active cell).
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$ws.Range("I44").Select()
